# Rename the two mapping sheets.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Complex Patient JSON to FHIR Ma")
$ws1.Name = "Patient-V1"
$ws2 = $wb.Worksheets.Item("Sheet1")
$ws2.Name = "Patient-V2"

# Make "Patient-V1" the active sheet/tab (was "Sheet1"/"Patient-V2"),
# and move its selection to A7.
$ws1.Activate()
$ws1.Range("A7").Select() | Out-Null
